$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "321.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "7.98%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "47.93"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "14.69%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.239"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.53%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08107"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "7.77%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.597"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.38%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.650"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.46%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.208"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "32.02%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1310"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "11.12%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1941"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.46%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09558"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "7.64%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "11.44%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1049"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.04%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001332"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "3.20%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005943"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.28%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.363"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.55%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.44%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3403"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.22%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.225"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.86%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1409"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04297"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.84%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001306"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.04%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004252"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "9.39%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001350"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.65%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003538"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-4.99%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02672"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "11.79%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05620"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.78%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006298"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-9.92%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007702"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.94%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1440"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.86%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007678"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.54%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.24%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.36%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006987"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.99%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.19%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "32.41%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003999"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-4.93%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.19%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.19%"
